$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.424.67'
$ws.Range("E2").Value = '  -5.35%  '
$ws.Range("D3").Value = '3.092.85'
$ws.Range("E3").Value = '  -8.10%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '''505.82'
$ws.Range("E5").Value = '  -3.95%  '
$ws.Range("D6").Value = '''165.73'
$ws.Range("E6").Value = '  -9.93%  '
$ws.Range("E7").Value = '  -3.40%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '3.091.51'
$ws.Range("E9").Value = '  -7.98%  '
$ws.Range("D10").Value = '''0.576'
$ws.Range("E10").Value = '  -7.00%  '
$ws.Range("D11").Value = '''50.78'
$ws.Range("E11").Value = '  -11.28%  '
$ws.Range("D12").Value = '''0.124'
$ws.Range("E12").Value = '  -5.49%  '
$ws.Range("E13").Value = '  -4.62%  '
$ws.Range("E14").Value = '  -5.88%  '
$ws.Range("D15").Value = '3.585.73'
$ws.Range("E15").Value = '  -8.39%  '
$ws.Range("E16").Value = '  -8.62%  '
$ws.Range("D17").Value = '3.092.64'
$ws.Range("E17").Value = '  -8.34%  '
$ws.Range("D18").Value = '61.187.71'
$ws.Range("E18").Value = '  -5.44%  '
$ws.Range("D19").Value = '''16.48'
$ws.Range("E19").Value = '  -4.38%  '
$ws.Range("D20").Value = '''10.49'
$ws.Range("E20").Value = '  -4.67%  '
$ws.Range("E21").Value = '  -3.85%  '
$ws.Range("D22").Value = '''352.46'
$ws.Range("E22").Value = '  -4.43%  '
$ws.Range("D23").Value = '''77.99'
$ws.Range("E23").Value = '  -3.24%  '
$ws.Range("E24").Value = '  -3.08%  '
$ws.Range("D25").Value = '''10.66'
$ws.Range("E25").Value = '  -1.10%  '
$ws.Range("D26").Value = '''6.09'
$ws.Range("E26").Value = '  +4.64%  '
$ws.Range("D27").Value = '''3.79'
$ws.Range("E27").Value = '  +1.55%  '
$ws.Range("E28").Value = '  -4.06%  '
$ws.Range("D29").Value = '''10.71'
$ws.Range("E29").Value = '  -5.67%  '
$ws.Range("E30").Value = '  -7.80%  '
$ws.Range("D31").Value = '''623.18'
$ws.Range("E31").Value = '  -5.86%  '
$ws.Range("D32").Value = '''27.11'
$ws.Range("E32").Value = '  -7.46%  '
$ws.Range("D33").Value = '''6.12'
$ws.Range("E33").Value = '  -7.75%  '
$ws.Range("E34").Value = '  -2.28%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").Value = '''0.0997'
$ws.Range("E36").Value = '  -4.30%  '
$ws.Range("D37").Value = '''55.57'
$ws.Range("E37").Value = '  -8.93%  '
$ws.Range("D38").Value = '''34.79'
$ws.Range("E38").Value = '  -4.03%  '
$ws.Range("E39").Value = '  -4.26%  '
$ws.Range("D40").Value = '''0.998'
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").Value = '0.0₃0656'
$ws.Range("E41").Value = '  +5.43%  '
$ws.Range("D42").Value = '''0.117'
$ws.Range("E42").Value = '  -6.92%  '
$ws.Range("D43").Value = '2.760.85'
$ws.Range("E43").Value = '  -2.30%  '
$ws.Range("D44").Value = '''2.41'
$ws.Range("E44").Value = '  +4.35%  '
$ws.Range("D45").Value = '''2.81'
$ws.Range("E45").Value = '  +9.73%  '
$ws.Range("D46").Value = '''2.58'
$ws.Range("E46").Value = '  -1.50%  '
$ws.Range("E47").Value = '  -3.89%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").Value = '''2.86'
$ws.Range("E48").Value = '  +0.12%  '
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").Value = '''2.43'
$ws.Range("E49").Value = '  -10.79%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = '''0.119'
$ws.Range("E50").Value = '  -3.88%  '
$ws.Range("D51").Value = '''128.96'
$ws.Range("E51").Value = '  -6.32%  '
